$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    "E2" = 3
    "F2" = 1
    "G2" = 1.980997
    "H2" = 5.942991000000001
    "I2" = 0.2896906247626733
    "J2" = 0.2896906247626732
    "M2" = 0.07195966666666666
    "N2" = 0.215879
    "O2" = 0.07530091904660251
    "P2" = 0.07530091904660252
    "Q2" = 0.1425518837876667
    "R2" = 1.282966954089
    "S2" = 0.02181397028381376
    "T2" = 0.02181397028381376
    "E3" = 3
    "F3" = 1
    "G3" = 1.980997
    "H3" = 5.942991000000001
    "I3" = 0.2896906247626733
    "J3" = 0.2896906247626732
    "O3" = 0.2743421080169271
    "P3" = 0.2743421080169271
    "Q3" = 0.5193560024930001
    "R3" = 4.674204022437
    "S3" = 0.0794743366701324
    "T3" = 0.07947433667013239
    "E4" = 3
    "F4" = 1
    "G4" = 1.980997
    "H4" = 5.942991000000001
    "I4" = 0.2896906247626733
    "J4" = 0.2896906247626732
    "K4" = 3
    "L4" = 1
    "M4" = 0.6214993333333333
    "N4" = 1.864498
    "O4" = 0.6503569729364704
    "P4" = 0.6503569729364704
    "Q4" = 1.231188314835334
    "R4" = 11.080694833518
    "S4" = 0.1884023178087271
    "T4" = 0.1884023178087271
    "I5" = 0.3000904023298512
    "J5" = 0.3000904023298512
    "M5" = 0.07195966666666666
    "N5" = 0.215879
    "O5" = 0.07530091904660251
    "P5" = 0.07530091904660252
    "Q5" = 0.147669439402
    "R5" = 1.329024954618
    "S5" = 0.02259708309250251
    "T5" = 0.02259708309250251
    "I6" = 0.3000904023298512
    "J6" = 0.3000904023298512
    "O6" = 0.2743421080169271
    "P6" = 0.2743421080169271
    "S6" = 0.08232743357081916
    "T6" = 0.08232743357081916
    "I7" = 0.3000904023298512
    "J7" = 0.3000904023298512
    "K7" = 3
    "L7" = 1
    "M7" = 0.6214993333333333
    "N7" = 1.864498
    "O7" = 0.6503569729364704
    "P7" = 0.6503569729364704
    "Q7" = 1.275387482924
    "R7" = 11.478487346316
    "S7" = 0.1951658856665296
    "T7" = 0.1951658856665296
    "G8" = 2.610033666666667
    "H8" = 7.830101
    "I8" = 0.3816776519844691
    "J8" = 0.3816776519844691
    "M8" = 0.07195966666666666
    "N8" = 0.215879
    "O8" = 0.07530091904660251
    "P8" = 0.07530091904660252
    "Q8" = 0.1878171526421111
    "R8" = 1.690354373779
    "S8" = 0.02874067797397983
    "T8" = 0.02874067797397984
    "G9" = 2.610033666666667
    "H9" = 7.830101
    "I9" = 0.3816776519844691
    "J9" = 0.3816776519844691
    "O9" = 0.2743421080169271
    "P9" = 0.2743421080169271
    "Q9" = 0.6842699163563334
    "R9" = 6.158429247207
    "S9" = 0.1047102516283703
    "T9" = 0.1047102516283703
    "G10" = 2.610033666666667
    "H10" = 7.830101
    "I10" = 0.3816776519844691
    "J10" = 0.3816776519844691
    "K10" = 3
    "L10" = 1
    "M10" = 0.6214993333333333
    "N10" = 1.864498
    "O10" = 0.6503569729364704
    "P10" = 0.6503569729364704
    "Q10" = 1.622134183810889
    "R10" = 14.599207654298
    "S10" = 0.248226722382119
    "T10" = 0.2482267223821189
    "E11" = 1
    "F11" = 0.3333333333333333
    "G11" = 0.1951746666666667
    "H11" = 0.585524
    "I11" = 0.02854132092300653
    "J11" = 0.02854132092300652
    "M11" = 0.07195966666666666
    "N11" = 0.215879
    "O11" = 0.07530091904660251
    "P11" = 0.07530091904660252
    "Q11" = 0.01404470395511111
    "R11" = 0.126402335596
    "S11" = 0.002149187696306417
    "T11" = 0.002149187696306417
    "E12" = 1
    "F12" = 0.3333333333333333
    "G12" = 0.1951746666666667
    "H12" = 0.585524
    "I12" = 0.02854132092300653
    "J12" = 0.02854132092300652
    "O12" = 0.2743421080169271
    "P12" = 0.2743421080169271
    "Q12" = 0.05116874718533334
    "R12" = 0.460518724668
    "S12" = 0.007830086147605239
    "T12" = 0.007830086147605237
    "E13" = 1
    "F13" = 0.3333333333333333
    "G13" = 0.1951746666666667
    "H13" = 0.585524
    "I13" = 0.02854132092300653
    "J13" = 0.02854132092300652
    "K13" = 3
    "L13" = 1
    "M13" = 0.6214993333333333
    "N13" = 1.864498
    "O13" = 0.6503569729364704
    "P13" = 0.6503569729364704
    "Q13" = 0.1213009252168889
    "R13" = 1.091708326952
    "S13" = 0.01856204707909487
    "T13" = 0.01856204707909487
}

foreach ($key in $values.Keys) {
    $ws.Range($key).Value = $values[$key]
}

$wb.Save()